$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix - [Item] StorageType과 ZoneType 통합 (merge StorageType + ZoneType)
# Rename the "StorageType" column header (row 3) to "ZoneType"
$ws.Range("B3").Value = "ZoneType"

# Every data row's former StorageType code is consolidated to a single
# ZoneType value of 1
$ws.Range("B5:B15").Value = 1

# Leftover selection artifact from the author's last click in Excel
$ws.Range("C19").Select()
